# Update the "Estado de Cuenta" worksheet:
#  - Refresh the totals header (Valor Mora / Cant. Trabajadores / Cant. Periodos)
#  - Replace the detail table: remove the previous period's rows and add the
#    new period's workers, plus a back-log of the still-outstanding prior
#    periods for the original worker (Dionis Altamiranda).
# This mirrors the commit "Elimna EC anteriores y se agregan nuevos, se
# modifica base de datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header totals block
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 916460      # VALOR MORA
$ws.Range("C13").Value = 13          # Cant. Trabajadores
$ws.Range("F13").Value = 6           # Cant. Periodos

# ---------------------------------------------------------------------
# 2. Make room for the new detail rows.
#    Before: 14 data rows (16-29), last one carrying the "closing" border
#    style. After: 18 data rows (16-33), same closing style kept on the
#    (now last) row. Insert 4 blank rows right before the closing row and
#    stamp them with the regular data-row formatting copied from row 28.
# ---------------------------------------------------------------------
$ws.Rows("29:32").Insert()
$ws.Range("B28:J28").Copy()
$ws.Range("B29:J32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Detail table values (columns B..G); H:J stay blank as before.
# ---------------------------------------------------------------------
$data = @(
    @("CC","34988844","YENIS ESTHER HERRERA CATALAN","2507",56940,1423500),
    @("CC","1051824956","ADRIANA MARCELA ARIAS MEJIA","2507",72000,877803),
    @("CC","64576906","AYDA PEREIRA OTERO","2507",52000,1300000),
    @("CC","1128050183","KELLY JOHANA RIVERA GRAU","2507",40000,1000000),
    @("CC","78744748","JUAN CARLOS HERRERA CATALAN","2507",56940,877803),
    @("CC","1143382753","MELISSA ANDREA DIAZ HERRERA","2507",56940,1423500),
    @("CC","1043651362","GABRIELA CAROLINA ARGARIN TRONCOSO","2507",56940,1423500),
    @("CC","6893110","RAMIRO JOSE ARROYO HERRERA","2507",56940,1423500),
    @("CC","45561034","YESICA LICETH CERVANTES SALCEDO","2507",56940,1423500),
    @("CC","32936496","KARINA DEL CARMEN PEÑA PEREZ","2507",56940,828116),
    @("CC","1047482952","DIONIS ALTAMIRANDA MANJARRES","2507",40000,1000000),
    @("CC","1047482952","DIONIS ALTAMIRANDA MANJARRES","2506",40000,1000000),
    @("CC","1047482952","DIONIS ALTAMIRANDA MANJARRES","2505",40000,1000000),
    @("CC","1047482952","DIONIS ALTAMIRANDA MANJARRES","2504",40000,1000000),
    @("CC","1047482952","DIONIS ALTAMIRANDA MANJARRES","2503",40000,1000000),
    @("CC","1047482952","DIONIS ALTAMIRANDA MANJARRES","2502",40000,1000000),
    @("CC","1003050853","YAMIL ANDRES HERRERA PEREZ","2507",56940,1423500),
    @("CC","1003050121","CAMILA ANDREA HERRERA RIVERO","2507",56940,1423500)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4. Re-measure the bestFit columns for the new (generally longer) content.
#    ColumnWidth is expressed in characters (no gridline padding); the
#    engine re-adds the standard 5px padding when it serialises the
#    worksheet <col> width, so back the padding out of the numbers first.
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 17.709635416666668
$ws.Columns("C").ColumnWidth = 15.893229166666666
$ws.Columns("E").ColumnWidth = 12.709635416666666
$ws.Columns("F").ColumnWidth = 9.346354166666666
$ws.Columns("G").ColumnWidth = 13.529947916666666
$ws.Columns("H").ColumnWidth = 18.529947916666668
$ws.Columns("I").ColumnWidth = 17.256510416666668
$ws.Columns("J").ColumnWidth = 14.166666666666666
